$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move the existing B2 hyperlink target ("http://172.31.21.30:8083") down to A3
$ws.Range("A3").Value = "http://172.31.21.30:8083"
$ws.Hyperlinks.Add($ws.Range("A3"), "http://172.31.21.30:8083/", [Type]::Missing, [Type]::Missing, "http://172.31.21.30:8083")
$ws.Range("A3").Style = "Hyperlink"

# Add the new rows for the second host
$ws.Range("A4").Value = "http://172.31.25.94:8082/webapp/"
$ws.Hyperlinks.Add($ws.Range("A4"), "http://172.31.25.94:8082/webapp/", [Type]::Missing, [Type]::Missing, "http://172.31.25.94:8082/webapp/")
$ws.Range("A4").Style = "Hyperlink"

$ws.Range("A5").Value = "http://172.31.25.94:8083"
$ws.Hyperlinks.Add($ws.Range("A5"), "http://172.31.25.94:8083/", [Type]::Missing, [Type]::Missing, "http://172.31.25.94:8083")
$ws.Range("A5").Style = "Hyperlink"

# Remove header B1 text and the now-relocated B2, then delete the now-empty column B
$ws.Columns("B").Delete()
